$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.025.37"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "1.823.54"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").Value = "'308.87"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").Value = "'0.4670"
$ws.Range("E7").Value = "  -1.58%  "

$ws.Range("D8").Value = "'0.3663"
$ws.Range("E8").Value = "  -0.64%  "

$ws.Range("D9").Value = "'0.07242"

$ws.Range("D10").Value = "'0.8608"
$ws.Range("E10").Value = "  -2.72%  "

$ws.Range("D11").Value = "'19.89"
$ws.Range("E11").Value = "  -2.66%  "

$ws.Range("D12").Value = "'0.07544"
$ws.Range("E12").Value = "  +2.92%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.337"
$ws.Range("E13").Value = "  -1.94%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'91.90"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.487"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.702.07"
$ws.Range("E16").Value = "  -9.02%  "

$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").Value = "'0.000008649"
$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'14.49"
$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.812.37"
$ws.Range("E21").Value = "  -2.63%  "

$ws.Range("D22").Value = "'5.152"
$ws.Range("E22").Value = "  -2.68%  "

$ws.Range("D23").Value = "'10.53"
$ws.Range("E23").Value = "  -1.23%  "

$ws.Range("D24").Value = "1.968.98"
$ws.Range("E24").Value = "  -5.85%  "

$ws.Range("D25").Value = "'151.39"
$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("D26").Value = "'1.840"
$ws.Range("E26").Value = "  -2.79%  "

$ws.Range("D27").Value = "'18.18"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("D28").Value = "'2.063"
$ws.Range("E28").Value = "  -3.56%  "

$ws.Range("D29").Value = "'5.116"
$ws.Range("E29").Value = "  -2.51%  "

$ws.Range("D30").Value = "'115.23"
$ws.Range("E30").Value = "  -1.86%  "

$ws.Range("D31").Value = "'0.08869"
$ws.Range("E31").Value = "  -1.41%  "

$ws.Range("D32").Value = "'2.956"
$ws.Range("E32").Value = "  +0.23%  "

$ws.Range("D33").Value = "'4.424"
$ws.Range("E33").Value = "  -2.74%  "

$ws.Range("D34").Value = "'1.134"
$ws.Range("E34").Value = "  -3.59%  "

$ws.Range("D35").Value = "'0.7200"
$ws.Range("E35").Value = "  -4.64%  "

$ws.Range("D36").Value = "'1.081"
$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("D37").Value = "'0.05261"
$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("D38").Value = "'2.414"
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").Value = "'0.01925"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("D40").Value = "'2.925"
$ws.Range("E40").Value = "  -1.67%  "

$ws.Range("D41").Value = "'7.164"
$ws.Range("E41").Value = "  -1.68%  "

$ws.Range("D42").Value = "'0.5165"
$ws.Range("E42").Value = "  -2.85%  "

$ws.Range("D43").Value = "'0.1632"
$ws.Range("E43").Value = "  -1.60%  "

$ws.Range("D44").Value = "'0.8580"
$ws.Range("E44").Value = "  -15.09%  "

$ws.Range("D45").Value = "'8.183"
$ws.Range("E45").Value = "  -3.46%  "

$ws.Range("D46").Value = "'0.4817"
$ws.Range("E46").Value = "  -1.96%  "

$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("D48").Value = "'10.13"
$ws.Range("E48").Value = "  -4.00%  "

$ws.Range("D49").Value = "'102.86"
$ws.Range("E49").Value = "  -2.01%  "

$ws.Range("E50").Value = "  -2.95%  "

$ws.Range("D51").Value = "'0.06240"
$ws.Range("E51").Value = "  -0.90%  "
